$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...would love to have you as majors: " ->
#           "...would love to have you among our majors: "
#
# Locate the exact character offsets of the word boundaries inside the
# sentence, drop temporary bookmarks at those boundaries (which splits the
# run there without merging neighboring runs back together), do the small
# text edit confined to a single split-off run, then remove the temporary
# bookmarks again. This reproduces the same run layout a live edit in Word
# would leave behind.
# ---------------------------------------------------------------------------
$needle1 = "you as majors: "
$find1 = $d.Content
$find1.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base1 = $find1.Start

$posAEnd     = $base1 + ("you a".Length)          # right after "a" of "as"
$posSEnd     = $posAEnd + 1                        # right after "s" of "as"
$posMajorEnd = $base1 + ("you as major".Length)    # right after "major"
$posSEnd2    = $posMajorEnd + 1                     # right after the trailing "s"

$d.Bookmarks.Add("zzTmp1", $d.Range($posAEnd, $posAEnd))
$d.Bookmarks.Add("zzTmp2", $d.Range($posSEnd, $posSEnd))
$d.Bookmarks.Add("zzTmp3", $d.Range($posMajorEnd, $posMajorEnd))
$d.Bookmarks.Add("zzTmp4", $d.Range($posSEnd2, $posSEnd2))

$sRange = $d.Range($posAEnd, $posSEnd)
$sRange.Text = "mong our"

$d.Bookmarks.Item("zzTmp1").Delete()
$d.Bookmarks.Item("zzTmp2").Delete()
$d.Bookmarks.Item("zzTmp3").Delete()
$d.Bookmarks.Item("zzTmp4").Delete()

# ---------------------------------------------------------------------------
# Change 2 & 3: move the "_GoBack" bookmark from the trailing empty paragraph
# into the middle of "majors." ("...hundreds of majo|rs.") in the paragraph
# about program size, leaving the trailing paragraph fully empty.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$needle2 = "hundreds of majors."
$find2 = $d.Content
$find2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $find2.Start + ($needle2.Length - 3)

$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
